$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Wipe the old sample content (keeps per-cell styles s="1"/s="3" in place,
#    so we don't disturb the style table while we repopulate).
# ---------------------------------------------------------------------------
$ws.Range("A1:F3").ClearContents()

# ---------------------------------------------------------------------------
# 2. Product attribute strings, written in (product, field) order so the
#    regenerated sharedStrings table lines up with the target ordering:
#    name, description, category, unit -- for products A, B, C, D.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "Sản phẩm A"
$ws.Range("C2").Value = "Mô tả sản phẩm A"
$ws.Range("F2").Value = "Danh mục A"
$ws.Range("G2").Value = "cái"

$ws.Range("B3").Value = "Sản phẩm B"
$ws.Range("C3").Value = "Mô tả sản phẩm B"
$ws.Range("F3").Value = "Danh mục B"
$ws.Range("G3").Value = "hộp"

$ws.Range("B4").Value = "Sản phẩm C"
$ws.Range("C4").Value = "Mô tả sản phẩm C"
$ws.Range("F4").Value = "Danh mục C"
$ws.Range("G4").Value = "kg"

$ws.Range("B5").Value = "Sản phẩm D"
$ws.Range("C5").Value = "Mô tả sản phẩm D"
$ws.Range("F5").Value = "Chưa phân loại"
$ws.Range("G5").Value = "gói"

# ---------------------------------------------------------------------------
# 3. Header row, written in the new header order.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("A1").Value = "productCode"
$ws.Range("C1").Value = "description"
$ws.Range("D1").Value = "price"
$ws.Range("E1").Value = "stock"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "unit"
$ws.Range("H1").Value = "discountPercentage"

# ---------------------------------------------------------------------------
# 4. Product codes (kept for last, matches the order new strings land at
#    the tail of the shared-string table).
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "SAN-2908"
$ws.Range("A5").Value = "SAN-9564"
$ws.Range("A3").Value = "SAN-5425"
$ws.Range("A2").Value = "SAN-2712"

# ---------------------------------------------------------------------------
# 5. Numeric columns: price, stock, discountPercentage.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 100000
$ws.Range("E2").Value = 50
$ws.Range("H2").Value = 10

$ws.Range("D3").Value = 200000
$ws.Range("E3").Value = 30
$ws.Range("H3").Value = 5

$ws.Range("D4").Value = 150000
$ws.Range("E4").Value = 20
$ws.Range("H4").Value = 15

$ws.Range("D5").Value = 50000
$ws.Range("E5").Value = 100
$ws.Range("H5").Value = 0

# ---------------------------------------------------------------------------
# 6. Formatting: header row bold/centered/wrapped (already s=1, but the row
#    grew from 6 to 8 columns and got taller), data rows keep the existing
#    "centered" look but drop the horizontal centering (Excel re-derives a
#    fresh style the first time an existing alignment combo changes), and
#    the fresh blank templating row (row 6) gets left/right aligned cells.
# ---------------------------------------------------------------------------
$ws.Range("A1:H1").Font.Bold = $true
$ws.Range("A1:H1").HorizontalAlignment = -4108
$ws.Range("A1:H1").VerticalAlignment = -4108
$ws.Range("A1:H1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 45

$ws.Range("A2:H5").VerticalAlignment = -4108
$ws.Range("A2:H5").WrapText = $true
$ws.Range("A2:H5").HorizontalAlignment = 1

$ws.Range("A6:H6").HorizontalAlignment = -4131
$ws.Range("D6").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 7. View / selection bookkeeping to match the edited workbook.
# ---------------------------------------------------------------------------
$ws.Range("C9").Select()

Write-Output "edit applied"
